$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 766, shifting the
# existing rows 766..794 down to 768..796 (matches the diff, which shows
# every existing row from 766 onward taking on the values previously held
# by the row two above its new position, plus two brand new rows of data
# and an extended used range A1:R796).
$ws.Rows("766:767").Insert()

# New row 766
$ws.Range("A766").Value = 10
$ws.Range("B766").Value = "Vega Modelo de Temuco"
$ws.Range("C766").Value = "La Araucanía"
$ws.Range("D766").Value = 45075
$ws.Range("E766").Value = 9
$ws.Range("F766").Value = 100112043
$ws.Range("G766").Value = "Pepino ensalada"
$ws.Range("H766").Value = "Alaska"
$ws.Range("I766").Value = "Primera"
$ws.Range("J766").Value = 35
$ws.Range("K766").Value = 28000
$ws.Range("L766").Value = 28000
$ws.Range("M766").Value = 28000
$ws.Range("N766").Value = "$/caja 60 unidades"
$ws.Range("O766").Value = "Región de Arica y Parinacota"
$ws.Range("P766").Value = 467
$ws.Range("Q766").Value = 60
$ws.Range("R766").Value = "Hortaliza"

# New row 767
$ws.Range("A767").Value = 10
$ws.Range("B767").Value = "Vega Modelo de Temuco"
$ws.Range("C767").Value = "La Araucanía"
$ws.Range("D767").Value = 45075
$ws.Range("E767").Value = 9
$ws.Range("F767").Value = 100112043
$ws.Range("G767").Value = "Pepino ensalada"
$ws.Range("H767").Value = "Sin especificar"
$ws.Range("I767").Value = "Primera"
$ws.Range("J767").Value = 345
$ws.Range("K767").Value = 13000
$ws.Range("L767").Value = 15000
$ws.Range("M767").Value = 14275
$ws.Range("N767").Value = "$/caja 60 unidades"
$ws.Range("O767").Value = "Región de Arica y Parinacota"
$ws.Range("P767").Value = 238
$ws.Range("Q767").Value = 60
$ws.Range("R767").Value = "Hortaliza"

# Apply the same date number format (yyyy-mm-dd hh:mm:ss) used by the
# rest of column D to the two newly inserted date cells.
$ws.Range("D766:D767").NumberFormat = $ws.Range("D765").NumberFormat
